$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.986.27'
$ws.Range("E2").Value = '  +0.76%  '

$ws.Range("D3").Value = '3.848.22'
$ws.Range("E3").Value = '  +1.40%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '689.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.27%  '

$ws.Range("D7").Value = '3.847.13'
$ws.Range("E7").Value = '  +1.44%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -0.08%  '

$ws.Range("E10").Value = '  +1.62%  '

$ws.Range("E11").Value = '  +5.19%  '

$ws.Range("E12").Value = '  -0.37%  '

$ws.Range("E13").Value = '  +6.00%  '

$ws.Range("E14").Value = '  +2.55%  '

$ws.Range("D15").Value = '4.490.33'
$ws.Range("E15").Value = '  +1.23%  '

$ws.Range("D16").Value = '3.843.96'
$ws.Range("E16").Value = '  +0.90%  '

$ws.Range("D17").Value = '71.014.50'
$ws.Range("E17").Value = '  +0.84%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.76'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.68%  '

$ws.Range("E19").Value = '  +0.70%  '

$ws.Range("E20").Value = '  +0.33%  '

$ws.Range("E21").Value = '  -3.70%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '489.82'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.721'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.66'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.25%  '

$ws.Range("E25").Value = '  +3.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.53'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.15'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.38%  '

$ws.Range("D29").Value = '4.000.67'
$ws.Range("E29").Value = '  +1.33%  '

$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("E31").Value = '  +9.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.63'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.91%  '

$ws.Range("E33").Value = '  +0.22%  '

$ws.Range("E34").Value = '  +0.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.181'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.27'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.98%  '

$ws.Range("D37").Value = '3.799.04'
$ws.Range("E37").Value = '  +1.20%  '

$ws.Range("E38").Value = '  +0.02%  '

$ws.Range("E39").Value = '  +1.52%  '

$ws.Range("E40").Value = '  +12.57%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.43'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.07'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.01%  '

$ws.Range("E43").Value = '  +5.08%  '

$ws.Range("E44").Value = '  -0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '165.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.05%  '

$ws.Range("E47").Value = '  +7.73%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.66'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.47%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.45%  '

$ws.Range("E50").Value = '  +1.26%  '

$ws.Range("E51").Value = '  -2.10%  '
